# OGD-Metadaten_Geschaeftsberichte.xlsx update
# - "Update for Prod and adding categories to metadata"
#
# 1. On the "metadata" sheet, the example/default "Kategorie" value (D4)
#    changes from the old free-text form "Politik, Verwaltung" to the new
#    machine-readable, comma-joined slug form "politik, verwaltung".
# 2. On the "kategorien_werteliste" sheet, a new "Code" column (B) is added
#    next to the existing category display names (A), giving each category
#    a slugified code value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) kategorien_werteliste sheet: add "Code" column with slugs
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("kategorien_werteliste")

# Row -> slug code, in the same order as the existing category rows (A2:A20)
$codes = @(
    "Code",
    "arbeit-und-erwerb",
    "basiskarten,",
    "bauen-und-wohnen",
    "bevolkerung",
    "bildung",
    "energie",
    "finanzen",
    "freizeit",
    "gesundheit",
    "kriminalitat",
    "kultur",
    "mobilitat",
    "politik",
    "preise",
    "soziales",
    "tourismus",
    "umwelt",
    "verwaltung",
    "volkswirtschaft"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 1
    $srcCell = $ws2.Cells.Item($row, 1)
    $dstCell = $ws2.Cells.Item($row, 2)
    $dstCell.Value = $codes[$i]

    # Mirror the formatting of the corresponding column A cell on the same
    # row so the new column fits visually with the existing table.
    $dstCell.Font.Bold = $srcCell.Font.Bold
    $dstCell.Font.Size = $srcCell.Font.Size
    $dstCell.Font.Name = $srcCell.Font.Name
    $dstCell.Font.Color = $srcCell.Font.Color
    if ($srcCell.Interior.ColorIndex -eq -4142) {
        $dstCell.Interior.ColorIndex = -4142
    } else {
        $dstCell.Interior.Color = $srcCell.Interior.Color
    }
    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $dstCell.VerticalAlignment = $srcCell.VerticalAlignment
}

# ---------------------------------------------------------------------
# 2) metadata sheet: update the Kategorie example value
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("metadata")
$ws1.Range("D4").Value = "politik, verwaltung"

Write-Output "Updated metadata!D4 and kategorien_werteliste!B1:B20"
